# Update the pandapower vm_pu.xlsx results sheet for the "Case_3_79" run
# with a 380 kV slack-bus setpoint: the slack bus voltage (column B) moves
# from 1.05 p.u. to 1.02 p.u. and every other bus's resulting per-unit
# voltage (columns C-F, I-N) is refreshed with the newly recomputed values
# for rows 2-25 (the 24 timesteps in the series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.030490150665257
$ws.Cells.Item(2, 4).Value = 1.032076269318367
$ws.Cells.Item(2, 5).Value = 1.039721047873346
$ws.Cells.Item(2, 6).Value = 1.049724111409681
$ws.Cells.Item(2, 9).Value = 1.030255399101766
$ws.Cells.Item(2, 10).Value = 1.035631108068005
$ws.Cells.Item(2, 11).Value = 1.03488254408658
$ws.Cells.Item(2, 12).Value = 1.042505431779909
$ws.Cells.Item(2, 13).Value = 1.052480373809751
$ws.Cells.Item(2, 14).Value = 1.015797953245887
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031675592545348
$ws.Cells.Item(3, 4).Value = 1.032911225731452
$ws.Cells.Item(3, 5).Value = 1.040848051644135
$ws.Cells.Item(3, 6).Value = 1.051106086446118
$ws.Cells.Item(3, 9).Value = 1.030432489720505
$ws.Cells.Item(3, 10).Value = 1.036456822177965
$ws.Cells.Item(3, 11).Value = 1.035526127037439
$ws.Cells.Item(3, 12).Value = 1.043441877737671
$ws.Cells.Item(3, 13).Value = 1.053673182614326
$ws.Cells.Item(3, 14).Value = 1.016077635605407
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.032442305104613
$ws.Cells.Item(4, 4).Value = 1.033450948351937
$ws.Cells.Item(4, 5).Value = 1.04157736297901
$ws.Cells.Item(4, 6).Value = 1.052000810570163
$ws.Cells.Item(4, 9).Value = 1.0305454780581
$ws.Cells.Item(4, 10).Value = 1.036990271242839
$ws.Cells.Item(4, 11).Value = 1.035941376478102
$ws.Cells.Item(4, 12).Value = 1.044047304442812
$ws.Cells.Item(4, 13).Value = 1.054444973114209
$ws.Cells.Item(4, 14).Value = 1.016258166794103
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03276454983512
$ws.Cells.Item(5, 4).Value = 1.033677716486172
$ws.Cells.Item(5, 5).Value = 1.041883982381012
$ws.Cells.Item(5, 6).Value = 1.05237707317128
$ws.Cells.Item(5, 9).Value = 1.030592595129303
$ws.Cells.Item(5, 10).Value = 1.037214332440553
$ws.Cells.Item(5, 11).Value = 1.036115662599569
$ws.Cells.Item(5, 12).Value = 1.044301703523165
$ws.Cells.Item(5, 13).Value = 1.054769427047745
$ws.Cells.Item(5, 14).Value = 1.016333956550886
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032818651437737
$ws.Cells.Item(6, 4).Value = 1.033715784168905
$ws.Cells.Item(6, 5).Value = 1.04193546616241
$ws.Cells.Item(6, 6).Value = 1.052440256456689
$ws.Cells.Item(6, 9).Value = 1.030600483822878
$ws.Cells.Item(6, 10).Value = 1.037251941558214
$ws.Cells.Item(6, 11).Value = 1.036144909325174
$ws.Cells.Item(6, 12).Value = 1.044344411114732
$ws.Cells.Item(6, 13).Value = 1.05482390394946
$ws.Cells.Item(6, 14).Value = 1.016346675802291
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03244661127117
$ws.Cells.Item(7, 4).Value = 1.033453978954086
$ws.Cells.Item(7, 5).Value = 1.041581459974169
$ws.Cells.Item(7, 6).Value = 1.052005837732729
$ws.Cells.Item(7, 9).Value = 1.030546109144423
$ws.Cells.Item(7, 10).Value = 1.03699326594717
$ws.Cells.Item(7, 11).Value = 1.035943706415942
$ws.Cells.Item(7, 12).Value = 1.044050704214701
$ws.Cells.Item(7, 13).Value = 1.054449308507797
$ws.Cells.Item(7, 14).Value = 1.01625917991498
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030890849470345
$ws.Cells.Item(8, 4).Value = 1.032358560884742
$ws.Cells.Item(8, 5).Value = 1.040101911596124
$ws.Cells.Item(8, 6).Value = 1.050191055498615
$ws.Cells.Item(8, 9).Value = 1.030315579285295
$ws.Cells.Item(8, 10).Value = 1.035910337248063
$ws.Cells.Item(8, 11).Value = 1.035100293063721
$ws.Cells.Item(8, 12).Value = 1.042822015812582
$ws.Cells.Item(8, 13).Value = 1.052883497879407
$ws.Cells.Item(8, 14).Value = 1.015892564885362
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028146646806978
$ws.Cells.Item(9, 4).Value = 1.030424062362268
$ws.Cells.Item(9, 5).Value = 1.037495186451022
$ws.Cells.Item(9, 6).Value = 1.046996854451491
$ws.Cells.Item(9, 9).Value = 1.029897089124129
$ws.Cells.Item(9, 10).Value = 1.033995572418792
$ws.Cells.Item(9, 11).Value = 1.033604937502058
$ws.Cells.Item(9, 12).Value = 1.040652894440398
$ws.Cells.Item(9, 13).Value = 1.05012397149613
$ws.Cells.Item(9, 14).Value = 1.015243145251388
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.026315201053901
$ws.Cells.Item(10, 4).Value = 1.029131521072317
$ws.Cells.Item(10, 5).Value = 1.035757562577315
$ws.Cells.Item(10, 6).Value = 1.044869709922306
$ws.Cells.Item(10, 9).Value = 1.02960984061626
$ws.Cells.Item(10, 10).Value = 1.032714618249435
$ws.Cells.Item(10, 11).Value = 1.0326018407536
$ws.Cells.Item(10, 12).Value = 1.039204027575725
$ws.Cells.Item(10, 13).Value = 1.048283914543442
$ws.Cells.Item(10, 14).Value = 1.014807897123373
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025521666764405
$ws.Cells.Item(11, 4).Value = 1.028571146679589
$ws.Cells.Item(11, 5).Value = 1.035005176092133
$ws.Cells.Item(11, 6).Value = 1.043949149852537
$ws.Cells.Item(11, 9).Value = 1.02948349877982
$ws.Cells.Item(11, 10).Value = 1.03215888052797
$ws.Cells.Item(11, 11).Value = 1.032166010490619
$ws.Cells.Item(11, 12).Value = 1.038575971547479
$ws.Cells.Item(11, 13).Value = 1.047487032595368
$ws.Cells.Item(11, 14).Value = 1.014618879692115
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.02522683414516
$ws.Cells.Item(12, 4).Value = 1.02836289341855
$ws.Cells.Item(12, 5).Value = 1.034725706690316
$ws.Cells.Item(12, 6).Value = 1.043607286095838
$ws.Cells.Item(12, 9).Value = 1.029436274884896
$ws.Cells.Item(12, 10).Value = 1.031952291620383
$ws.Cells.Item(12, 11).Value = 1.032003900051467
$ws.Cells.Item(12, 12).Value = 1.038342578538069
$ws.Cells.Item(12, 13).Value = 1.047191013804892
$ws.Cells.Item(12, 14).Value = 1.014548586845552
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02529008034386
$ws.Cells.Item(13, 4).Value = 1.028407569225745
$ws.Cells.Item(13, 5).Value = 1.034785653837752
$ws.Cells.Item(13, 6).Value = 1.043680613766272
$ws.Cells.Item(13, 9).Value = 1.029446417919955
$ws.Cells.Item(13, 10).Value = 1.031996613050448
$ws.Cells.Item(13, 11).Value = 1.032038683427331
$ws.Cells.Item(13, 12).Value = 1.03839264691303
$ws.Cells.Item(13, 13).Value = 1.047254511847671
$ws.Cells.Item(13, 14).Value = 1.014563668678382
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025497297427681
$ws.Cells.Item(14, 4).Value = 1.028553934551856
$ws.Cells.Item(14, 5).Value = 1.034982075066743
$ws.Cells.Item(14, 6).Value = 1.043920889811117
$ws.Cells.Item(14, 9).Value = 1.029479601255452
$ws.Cells.Item(14, 10).Value = 1.032141807161816
$ws.Cells.Item(14, 11).Value = 1.032152614959306
$ws.Cells.Item(14, 12).Value = 1.038556681361688
$ws.Cells.Item(14, 13).Value = 1.047462564040202
$ws.Cells.Item(14, 14).Value = 1.014613070966303
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025624960211923
$ws.Cells.Item(15, 4).Value = 1.028644101088599
$ws.Cells.Item(15, 5).Value = 1.035103096667347
$ws.Cells.Item(15, 6).Value = 1.044068941427383
$ws.Cells.Item(15, 9).Value = 1.029500007517699
$ws.Cells.Item(15, 10).Value = 1.032231244378774
$ws.Cells.Item(15, 11).Value = 1.032222782270115
$ws.Cells.Item(15, 12).Value = 1.038657734426886
$ws.Cells.Item(15, 13).Value = 1.047590748956562
$ws.Cells.Item(15, 14).Value = 1.014643498288593
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.026367854424411
$ws.Cells.Item(16, 4).Value = 1.029168696556438
$ws.Cells.Item(16, 5).Value = 1.035807496160936
$ws.Cells.Item(16, 6).Value = 1.044930814898869
$ws.Cells.Item(16, 9).Value = 1.029618184164067
$ws.Cells.Item(16, 10).Value = 1.032751477931404
$ws.Cells.Item(16, 11).Value = 1.032630734055608
$ws.Cells.Item(16, 12).Value = 1.039245694955223
$ws.Cells.Item(16, 13).Value = 1.048336798119526
$ws.Cells.Item(16, 14).Value = 1.014820429920588
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.026833714979407
$ws.Cells.Item(17, 4).Value = 1.02949757468987
$ws.Cells.Item(17, 5).Value = 1.036249350270608
$ws.Cells.Item(17, 6).Value = 1.045471579117424
$ws.Cells.Item(17, 9).Value = 1.029691787975767
$ws.Cells.Item(17, 10).Value = 1.033077517531398
$ws.Cells.Item(17, 11).Value = 1.032886233819088
$ws.Cells.Item(17, 12).Value = 1.039614321828006
$ws.Cells.Item(17, 13).Value = 1.048804739744061
$ws.Cells.Item(17, 14).Value = 1.014931266354317
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02710539489647
$ws.Cells.Item(18, 4).Value = 1.029689336690514
$ws.Cells.Item(18, 5).Value = 1.036507077987905
$ws.Cells.Item(18, 6).Value = 1.04578704691322
$ws.Cells.Item(18, 9).Value = 1.029734530610277
$ws.Cells.Item(18, 10).Value = 1.033267587187261
$ws.Cells.Item(18, 11).Value = 1.033035119506988
$ws.Cells.Item(18, 12).Value = 1.039829269472904
$ws.Cells.Item(18, 13).Value = 1.049077670275902
$ws.Cells.Item(18, 14).Value = 1.014995862127404
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027198022551414
$ws.Cells.Item(19, 4).Value = 1.029754711227029
$ws.Cells.Item(19, 5).Value = 1.036594956825406
$ws.Cells.Item(19, 6).Value = 1.045894621710811
$ws.Cells.Item(19, 9).Value = 1.029749072655601
$ws.Cells.Item(19, 10).Value = 1.033332378479104
$ws.Cells.Item(19, 11).Value = 1.033085861450475
$ws.Cells.Item(19, 12).Value = 1.039902549895027
$ws.Cells.Item(19, 13).Value = 1.049170730592166
$ws.Cells.Item(19, 14).Value = 1.015017878597809
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.026783737604168
$ws.Cells.Item(20, 4).Value = 1.029462296120338
$ws.Cells.Item(20, 5).Value = 1.036201943346834
$ws.Cells.Item(20, 6).Value = 1.0454135551711
$ws.Cells.Item(20, 9).Value = 1.029683910557694
$ws.Cells.Item(20, 10).Value = 1.033042547308576
$ws.Cells.Item(20, 11).Value = 1.032858835911069
$ws.Cells.Item(20, 12).Value = 1.039574778511831
$ws.Cells.Item(20, 13).Value = 1.048754535286895
$ws.Cells.Item(20, 14).Value = 1.014919380164387
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025436279353364
$ws.Cells.Item(21, 4).Value = 1.028510836523337
$ws.Cells.Item(21, 5).Value = 1.034924233914324
$ws.Cells.Item(21, 6).Value = 1.043850132507048
$ws.Cells.Item(21, 9).Value = 1.029469837733185
$ws.Cells.Item(21, 10).Value = 1.032099055631354
$ws.Cells.Item(21, 11).Value = 1.03211907114463
$ws.Cells.Item(21, 12).Value = 1.038508380219055
$ws.Cells.Item(21, 13).Value = 1.047401298465067
$ws.Cells.Item(21, 14).Value = 1.014598525528163
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024588622520228
$ws.Cells.Item(22, 4).Value = 1.027912006715534
$ws.Cells.Item(22, 5).Value = 1.03412088883345
$ws.Cells.Item(22, 6).Value = 1.042867568879683
$ws.Cells.Item(22, 9).Value = 1.029333535262801
$ws.Cells.Item(22, 10).Value = 1.03150490014706
$ws.Cells.Item(22, 11).Value = 1.031652657110945
$ws.Cells.Item(22, 12).Value = 1.037837284680817
$ws.Cells.Item(22, 13).Value = 1.046550339950089
$ws.Cells.Item(22, 14).Value = 1.014396309004487
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025038025573721
$ws.Cells.Item(23, 4).Value = 1.028229515750875
$ws.Cells.Item(23, 5).Value = 1.03454675777789
$ws.Cells.Item(23, 6).Value = 1.043388405432755
$ws.Cells.Item(23, 9).Value = 1.029405953609069
$ws.Cells.Item(23, 10).Value = 1.031819963271889
$ws.Cells.Item(23, 11).Value = 1.031900035029723
$ws.Cells.Item(23, 12).Value = 1.03819310341677
$ws.Cells.Item(23, 13).Value = 1.047001461844992
$ws.Cells.Item(23, 14).Value = 1.01450355367993
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026806320379143
$ws.Cells.Item(24, 4).Value = 1.029478237194662
$ws.Cells.Item(24, 5).Value = 1.036223364489903
$ws.Cells.Item(24, 6).Value = 1.045439773542865
$ws.Cells.Item(24, 9).Value = 1.029687470608559
$ws.Cells.Item(24, 10).Value = 1.033058349167535
$ws.Cells.Item(24, 11).Value = 1.032871216288298
$ws.Cells.Item(24, 12).Value = 1.039592646638782
$ws.Cells.Item(24, 13).Value = 1.04877722055683
$ws.Cells.Item(24, 14).Value = 1.014924751186526
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.028856428135628
$ws.Cells.Item(25, 4).Value = 1.030924681192888
$ws.Cells.Item(25, 5).Value = 1.038169046402508
$ws.Cells.Item(25, 6).Value = 1.047822211597335
$ws.Cells.Item(25, 9).Value = 1.030006733323757
$ws.Cells.Item(25, 10).Value = 1.034491363209332
$ws.Cells.Item(25, 11).Value = 1.033992611458355
$ws.Cells.Item(25, 12).Value = 1.041214149012088
$ws.Cells.Item(25, 14).Value = 1.015411440177789

Write-Host "Updated res_bus vm_pu values for 380 kV case (rows 2-25)."
